$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.510.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").Value = "'2.947.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'588.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "'146.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.36%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'2.946.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.61%  "
$ws.Range("E9").Value = "  +3.38%  "
$ws.Range("D10").Value = "'7.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.32%  "
$ws.Range("D11").Value = "'0.150"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.76%  "
$ws.Range("D12").Value = "'0.435"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.70%  "
$ws.Range("D13").Value = "'0.0000233"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.86%  "
$ws.Range("D14").Value = "'32.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "'3.432.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "'62.464.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'6.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "'2.944.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("D20").Value = "'434.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.47%  "
$ws.Range("D21").Value = "'13.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").Value = "'0.662"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("D23").Value = "'6.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").Value = "'11.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.17%  "
$ws.Range("D25").Value = "'80.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").Value = "'11.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.03%  "
$ws.Range("E27").Value = "  +2.76%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'7.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.86%  "
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("E31").Value = "  +5.08%  "
$ws.Range("E32").Value = "  +18.92%  "
$ws.Range("E33").Value = "  +4.48%  "
$ws.Range("D34").Value = "'26.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.04%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "'0.988"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.59%  "
$ws.Range("D37").Value = "'5.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("E38").Value = "  +7.94%  "
$ws.Range("D39").Value = "'49.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("E40").Value = "  +6.51%  "
$ws.Range("D41").Value = "'8.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("E43").Value = "  +4.57%  "
$ws.Range("D44").Value = "'38.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "'135.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.18%  "
$ws.Range("D46").Value = "'2.689.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("D48").Value = "'355.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.16%  "
$ws.Range("E50").Value = "  +2.43%  "
$ws.Range("D51").Value = "'22.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.90%  "
